$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = 'D2','E2','D3','E3','E4','D5','E5','D6','E6','E7','E8','D9','E9','D10','E10','D11','E11','E12','E13','D14','E14','D15','E15','D16','E16','E17','D18','E18','D19','E19','D20','E20','D21','E21','D22','E22','D23','E23','D24','E24','D25','E25','E26','E27','D28','E28','E29','D30','E30','D31','E31','D32','E32','B33','C33','D33','E33','B34','C34','D34','E34','D35','E35','E36','D37','E37','D38','E38','E39','D40','E40','D41','E41','D42','E42','D43','E43','E44','D45','E45','D46','E46','D47','E47','D48','E48','E49','D50','E50','D51','E51'
foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '56.882.55'
$ws.Range('E2').Value = '  +3.75%  '
$ws.Range('D3').Value = '2.495.17'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '491.94'
$ws.Range('E5').Value = '  +2.41%  '
$ws.Range('D6').Value = '152.50'
$ws.Range('E6').Value = '  +10.61%  '
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('D9').Value = '2.507.92'
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('D10').Value = '5.73'
$ws.Range('E10').Value = '  +4.99%  '
$ws.Range('D11').Value = '0.0990'
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').Value = '2.932.92'
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('D15').Value = '56.851.95'
$ws.Range('E15').Value = '  +3.30%  '
$ws.Range('D16').Value = '21.21'
$ws.Range('E16').Value = '  +3.63%  '
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').Value = '2.506.67'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').Value = '4.54'
$ws.Range('E19').Value = '  +4.55%  '
$ws.Range('D20').Value = '10.28'
$ws.Range('E20').Value = '  +3.89%  '
$ws.Range('D21').Value = '320.70'
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').Value = '5.89'
$ws.Range('E23').Value = '  +3.81%  '
$ws.Range('D24').Value = '58.35'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('D25').Value = '0.409'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('E26').Value = '  -0.44%  '
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('D28').Value = '2.609.37'
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  +4.23%  '
$ws.Range('D30').Value = '0.0₃0809'
$ws.Range('E30').Value = '  +4.45%  '
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').Value = '151.87'
$ws.Range('E32').Value = '  +2.46%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '18.24'
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '1.52'
$ws.Range('E34').Value = '  +3.24%  '
$ws.Range('D35').Value = '5.28'
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('E36').Value = '  +4.99%  '
$ws.Range('D37').Value = '3.80'
$ws.Range('E37').Value = '  +4.51%  '
$ws.Range('D38').Value = '0.876'
$ws.Range('E38').Value = '  +2.59%  '
$ws.Range('E39').Value = '  +7.95%  '
$ws.Range('D40').Value = '34.22'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').Value = '3.51'
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('D42').Value = '0.0561'
$ws.Range('E42').Value = '  +2.96%  '
$ws.Range('D43').Value = '0.615'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').Value = '267.12'
$ws.Range('E45').Value = '  +6.07%  '
$ws.Range('D46').Value = '4.83'
$ws.Range('E46').Value = '  +4.41%  '
$ws.Range('D47').Value = '0.0937'
$ws.Range('E47').Value = '  +3.99%  '
$ws.Range('D48').Value = '10.21'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('E49').Value = '  +3.54%  '
$ws.Range('D50').Value = '17.89'
$ws.Range('E50').Value = '  +4.16%  '
$ws.Range('D51').Value = '1.891.14'
$ws.Range('E51').Value = '  -3.63%  '
